$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds crypto prices stored as plain text (e.g. "63.505.82",
# "1.00", "0.0₃0509") -- force text format before assigning so Excel
# does not reinterpret these as numbers/dates.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '63.505.82'
$ws.Range('E2').Value = '  +3.60%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.069.10'
$ws.Range('E3').Value = '  +2.91%  '

$ws.Range('E4').Value = '  -0.02%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '550.65'
$ws.Range('E5').Value = '  +2.92%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.71'
$ws.Range('E6').Value = '  +7.73%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.999'
$ws.Range('E7').Value = '  -0.07%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.065.39'
$ws.Range('E8').Value = '  +2.94%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.502'
$ws.Range('E9').Value = '  +1.58%  '

$ws.Range('E10').Value = '  +6.14%  '

$ws.Range('E11').Value = '  +3.34%  '

$ws.Range('E12').Value = '  +3.18%  '

$ws.Range('E13').Value = '  +3.09%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.87'
$ws.Range('E14').Value = '  +3.95%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.567.74'
$ws.Range('E15').Value = '  +2.59%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '63.508.31'
$ws.Range('E16').Value = '  +3.47%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '3.071.16'
$ws.Range('E17').Value = '  +2.58%  '

$ws.Range('E18').Value = '  -0.50%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.78'
$ws.Range('E19').Value = '  +3.07%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '483.75'
$ws.Range('E20').Value = '  +4.16%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.93'
$ws.Range('E21').Value = '  +6.05%  '

$ws.Range('E22').Value = '  +1.58%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.31'
$ws.Range('E23').Value = '  +6.39%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '81.01'
$ws.Range('E24').Value = '  +0.92%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.89'
$ws.Range('E25').Value = '  +8.61%  '

$ws.Range('E26').Value = '  +0.02%  '

$ws.Range('E27').Value = '  +4.08%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.94'
$ws.Range('E28').Value = '  +3.40%  '

$ws.Range('E30').Value = '  -0.04%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.19'
$ws.Range('E31').Value = '  +3.01%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.17'
$ws.Range('E32').Value = '  +1.70%  '

$ws.Range('E33').Value = '  +8.49%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.70'
$ws.Range('E34').Value = '  +5.59%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '55.40'
$ws.Range('E35').Value = '  +0.69%  '

$ws.Range('E36').Value = '  +2.98%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '467.48'
$ws.Range('E37').Value = '  +3.83%  '

$ws.Range('E38').Value = '  +5.37%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.0397'
$ws.Range('E39').Value = '  +4.14%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.024.08'
$ws.Range('E40').Value = '  -4.00%  '

$ws.Range('E41').Value = '  +0.12%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.22'
$ws.Range('E42').Value = '  +2.05%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.57'
$ws.Range('E43').Value = '  +6.22%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '27.69'
$ws.Range('E44').Value = '  +5.67%  '

$ws.Range('E45').Value = '  +6.14%  '

$ws.Range('E47').Value = '  +3.61%  '

$ws.Range('E48').Value = '  +2.90%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.51'
$ws.Range('E49').Value = '  -1.18%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0₃0509'
$ws.Range('E50').Value = '  +4.00%  '

$ws.Range('E51').Value = '  +4.87%  '

